$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lotusBloom = "('Lotus Bloom', ['Artifact', 'Suspend 3" + [char]0x2014 + "{0} (Rather than cast this card from your hand, pay {0} and exile it with three time counters on it. At the beginning of your upkeep, remove a time counter. When the last is removed, cast it without paying its mana cost.)', '{T}, Sacrifice Lotus Bloom: Add three mana of any one color.'])"

$suddenShock = "('Sudden Shock', ['{1}{R}', 'Instant', 'Split second (As long as this spell is on the stack, players can" + [char]0x2019 + "t cast spells or activate abilities that aren" + [char]0x2019 + "t mana abilities.)', 'Sudden Shock deals 2 damage to any target.'])"

# Delete rows 4 through 10 (shifts rows up automatically), leaving A1:A3
$ws.Range("A4:A10").EntireRow.Delete()

# Update A2 and A3 with the combined tuple-style text
$ws.Range("A2").Value = $lotusBloom
$ws.Range("A3").Value = $suddenShock
